# Atualizacao de bases das ligas, do dia: 17-03-2024 as 10:24
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 235-239 (re-sorted / refreshed odds data) ---
# Row 235
$ws.Range("A235").Value = 233
$ws.Range("B235").Value = 6861095
$ws.Range("C235").Value = "Romania Liga I"
$ws.Range("D235").Value = "Romania Liga I"
$ws.Range("E235").Value = 45359.625
$ws.Range("F235").Value = "FC Botosani"
$ws.Range("G235").Value = "Farul Constanta"
$ws.Range("H235").Value = 0
$ws.Range("I235").Value = 0
$ws.Range("J235").Value = "D"
$ws.Range("K235").Value = 3.75
$ws.Range("L235").Value = 3.4
$ws.Range("M235").Value = 1.909
$ws.Range("N235").Value = 3.1
$ws.Range("O235").Value = 3
$ws.Range("P235").Value = 2.375
$ws.Range("Q235").Value = 0.25
$ws.Range("R235").Value = 1.775
$ws.Range("S235").Value = 2.1
$ws.Range("T235").Value = 2
$ws.Range("U235").Value = 1.8
$ws.Range("V235").Value = 2.05
$ws.Range("W235").Value = -1
$ws.Range("X235").Value = 2
$ws.Range("Y235").Value = -1
$ws.Range("Z235").Value = 0.3875
$ws.Range("AA235").Value = -0.5
$ws.Range("AB235").Value = -1
$ws.Range("AC235").Value = 1.05

# Row 236
$ws.Range("A236").Value = 234
$ws.Range("B236").Value = 6852370
$ws.Range("C236").Value = "Romania Liga I"
$ws.Range("D236").Value = "Romania Liga I"
$ws.Range("E236").Value = 45359.625
$ws.Range("F236").Value = "Dinamo Bucharest"
$ws.Range("G236").Value = "ACS UTA Batrana Doamna"
$ws.Range("H236").Value = 1
$ws.Range("I236").Value = 0
$ws.Range("J236").Value = "H"
$ws.Range("K236").Value = 2.55
$ws.Range("L236").Value = 2.875
$ws.Range("M236").Value = 3
$ws.Range("N236").Value = 2.375
$ws.Range("O236").Value = 3
$ws.Range("P236").Value = 3.1
$ws.Range("Q236").Value = -0.25
$ws.Range("R236").Value = 2
$ws.Range("S236").Value = 1.85
$ws.Range("T236").Value = 2.25
$ws.Range("U236").Value = 1.975
$ws.Range("V236").Value = 1.875
$ws.Range("W236").Value = 1.375
$ws.Range("X236").Value = -1
$ws.Range("Y236").Value = -1
$ws.Range("Z236").Value = 1
$ws.Range("AA236").Value = -1
$ws.Range("AB236").Value = -1
$ws.Range("AC236").Value = 0.875

# Row 237
$ws.Range("A237").Value = 235
$ws.Range("B237").Value = 6865915
$ws.Range("C237").Value = "Romania Liga I"
$ws.Range("D237").Value = "Romania Liga I"
$ws.Range("E237").Value = 45359.625
$ws.Range("F237").Value = "FC Voluntari"
$ws.Range("G237").Value = "Universitatea Cluj"
$ws.Range("H237").Value = 0
$ws.Range("I237").Value = 0
$ws.Range("J237").Value = "D"
$ws.Range("K237").Value = 3.5
$ws.Range("L237").Value = 3.25
$ws.Range("M237").Value = 2.05
$ws.Range("N237").Value = 3.4
$ws.Range("O237").Value = 3.1
$ws.Range("P237").Value = 2.15
$ws.Range("Q237").Value = 0.25
$ws.Range("R237").Value = 1.975
$ws.Range("S237").Value = 1.875
$ws.Range("T237").Value = 2.25
$ws.Range("U237").Value = 2.05
$ws.Range("V237").Value = 1.75
$ws.Range("W237").Value = -1
$ws.Range("X237").Value = 2.1
$ws.Range("Y237").Value = -1
$ws.Range("Z237").Value = 0.4875
$ws.Range("AA237").Value = -0.5
$ws.Range("AB237").Value = -1
$ws.Range("AC237").Value = 0.75

# Row 238
$ws.Range("A238").Value = 236
$ws.Range("B238").Value = 6870268
$ws.Range("C238").Value = "Romania Liga I"
$ws.Range("D238").Value = "Romania Liga I"
$ws.Range("E238").Value = 45359.625
$ws.Range("F238").Value = "Petrolul Ploiesti"
$ws.Range("G238").Value = "ACS Sepsi"
$ws.Range("H238").Value = 1
$ws.Range("I238").Value = 2
$ws.Range("J238").Value = "A"
$ws.Range("K238").Value = 2.8
$ws.Range("L238").Value = 3
$ws.Range("M238").Value = 2.55
$ws.Range("N238").Value = 3
$ws.Range("O238").Value = 3.2
$ws.Range("P238").Value = 2.3
$ws.Range("Q238").Value = 0.25
$ws.Range("R238").Value = 1.85
$ws.Range("S238").Value = 2
$ws.Range("T238").Value = 2.25
$ws.Range("U238").Value = 1.875
$ws.Range("V238").Value = 1.975
$ws.Range("W238").Value = -1
$ws.Range("X238").Value = -1
$ws.Range("Y238").Value = 1.3
$ws.Range("Z238").Value = -1
$ws.Range("AA238").Value = 1
$ws.Range("AB238").Value = 0.875
$ws.Range("AC238").Value = -1

# Row 239
$ws.Range("A239").Value = 237
$ws.Range("B239").Value = 6836277
$ws.Range("C239").Value = "Romania Liga I"
$ws.Range("D239").Value = "Romania Liga I"
$ws.Range("E239").Value = 45359.625
$ws.Range("F239").Value = "CFR Cluj"
$ws.Range("G239").Value = "AFC Hermannstadt"
$ws.Range("H239").Value = 1
$ws.Range("I239").Value = 0
$ws.Range("J239").Value = "H"
$ws.Range("K239").Value = 1.7
$ws.Range("L239").Value = 3.4
$ws.Range("M239").Value = 5
$ws.Range("N239").Value = 1.65
$ws.Range("O239").Value = 3.5
$ws.Range("P239").Value = 5.25
$ws.Range("Q239").Value = -0.75
$ws.Range("R239").Value = 1.85
$ws.Range("S239").Value = 2
$ws.Range("T239").Value = 2.25
$ws.Range("U239").Value = 1.875
$ws.Range("V239").Value = 1.975
$ws.Range("W239").Value = 0.6499999999999999
$ws.Range("X239").Value = -1
$ws.Range("Y239").Value = -1
$ws.Range("Z239").Value = 0.425
$ws.Range("AA239").Value = -0.5
$ws.Range("AB239").Value = -1
$ws.Range("AC239").Value = 0.9750000000000001

# --- Append new rows 243-250 ---
# Row 243
$ws.Range("A243").Value = 241
$ws.Range("B243").Value = 7951557
$ws.Range("C243").Value = "Romania Liga I"
$ws.Range("D243").Value = "Romania Liga I"
$ws.Range("E243").Value = 45366.52083333334
$ws.Range("F243").Value = "Universitatea Cluj"
$ws.Range("G243").Value = "FC Botosani"
$ws.Range("H243").Value = 3
$ws.Range("I243").Value = 0
$ws.Range("J243").Value = "H"
$ws.Range("K243").Value = 1.615
$ws.Range("L243").Value = 3.6
$ws.Range("M243").Value = 5
$ws.Range("N243").Value = 1.909
$ws.Range("O243").Value = 2.75
$ws.Range("P243").Value = 5.25
$ws.Range("Q243").Value = -0.5
$ws.Range("R243").Value = 1.975
$ws.Range("S243").Value = 1.875
$ws.Range("T243").Value = 1.75
$ws.Range("U243").Value = 1.775
$ws.Range("V243").Value = 2.1
$ws.Range("W243").Value = 0.909
$ws.Range("X243").Value = -1
$ws.Range("Y243").Value = -1
$ws.Range("Z243").Value = 0.9750000000000001
$ws.Range("AA243").Value = -1
$ws.Range("AB243").Value = 0.7749999999999999
$ws.Range("AC243").Value = -1
$ws.Range("A235").Copy()
$ws.Range("A243").PasteSpecial(-4122)
$ws.Range("E235").Copy()
$ws.Range("E243").PasteSpecial(-4122)

# Row 244
$ws.Range("A244").Value = 242
$ws.Range("B244").Value = 7949044
$ws.Range("C244").Value = "Romania Liga I"
$ws.Range("D244").Value = "Romania Liga I"
$ws.Range("E244").Value = 45366.64583333334
$ws.Range("F244").Value = "Rapid Bucuresti"
$ws.Range("G244").Value = "Farul Constanta"
$ws.Range("H244").Value = 1
$ws.Range("I244").Value = 2
$ws.Range("J244").Value = "A"
$ws.Range("K244").Value = 1.75
$ws.Range("L244").Value = 3.5
$ws.Range("M244").Value = 4.2
$ws.Range("N244").Value = 1.6
$ws.Range("O244").Value = 3.75
$ws.Range("P244").Value = 5
$ws.Range("Q244").Value = -0.75
$ws.Range("R244").Value = 1.8
$ws.Range("S244").Value = 2.05
$ws.Range("T244").Value = 2.5
$ws.Range("U244").Value = 1.925
$ws.Range("V244").Value = 1.925
$ws.Range("W244").Value = -1
$ws.Range("X244").Value = -1
$ws.Range("Y244").Value = 4
$ws.Range("Z244").Value = -1
$ws.Range("AA244").Value = 1.05
$ws.Range("AB244").Value = 0.925
$ws.Range("AC244").Value = -1
$ws.Range("A235").Copy()
$ws.Range("A244").PasteSpecial(-4122)
$ws.Range("E235").Copy()
$ws.Range("E244").PasteSpecial(-4122)

# Row 245
$ws.Range("A245").Value = 243
$ws.Range("B245").Value = 7951558
$ws.Range("C245").Value = "Romania Liga I"
$ws.Range("D245").Value = "Romania Liga I"
$ws.Range("E245").Value = 45367.5
$ws.Range("F245").Value = "Otelul Galati"
$ws.Range("G245").Value = "CSM Politehnica Iasi"
$ws.Range("H245").Value = 1
$ws.Range("I245").Value = 0
$ws.Range("J245").Value = "H"
$ws.Range("K245").Value = 2.1
$ws.Range("L245").Value = 3.1
$ws.Range("M245").Value = 3.4
$ws.Range("N245").Value = 2
$ws.Range("O245").Value = 3
$ws.Range("P245").Value = 3.75
$ws.Range("Q245").Value = -0.25
$ws.Range("R245").Value = 1.75
$ws.Range("S245").Value = 2.05
$ws.Range("T245").Value = 2
$ws.Range("U245").Value = 2.025
$ws.Range("V245").Value = 1.825
$ws.Range("W245").Value = 1
$ws.Range("X245").Value = -1
$ws.Range("Y245").Value = -1
$ws.Range("Z245").Value = 0.75
$ws.Range("AA245").Value = -1
$ws.Range("AB245").Value = -1
$ws.Range("AC245").Value = 0.825
$ws.Range("A235").Copy()
$ws.Range("A245").PasteSpecial(-4122)
$ws.Range("E235").Copy()
$ws.Range("E245").PasteSpecial(-4122)

# Row 246
$ws.Range("A246").Value = 244
$ws.Range("B246").Value = 7949046
$ws.Range("C246").Value = "Romania Liga I"
$ws.Range("D246").Value = "Romania Liga I"
$ws.Range("E246").Value = 45367.625
$ws.Range("F246").Value = "CFR Cluj"
$ws.Range("G246").Value = "CS U Craiova"
$ws.Range("H246").Value = 1
$ws.Range("I246").Value = 2
$ws.Range("J246").Value = "A"
$ws.Range("K246").Value = 2.05
$ws.Range("L246").Value = 3.2
$ws.Range("M246").Value = 3.4
$ws.Range("N246").Value = 1.833
$ws.Range("O246").Value = 3.4
$ws.Range("P246").Value = 4
$ws.Range("Q246").Value = -0.5
$ws.Range("R246").Value = 1.9
$ws.Range("S246").Value = 1.95
$ws.Range("T246").Value = 2.5
$ws.Range("U246").Value = 1.95
$ws.Range("V246").Value = 1.9
$ws.Range("W246").Value = -1
$ws.Range("X246").Value = -1
$ws.Range("Y246").Value = 3
$ws.Range("Z246").Value = -1
$ws.Range("AA246").Value = 0.95
$ws.Range("AB246").Value = 0.95
$ws.Range("AC246").Value = -1
$ws.Range("A235").Copy()
$ws.Range("A246").PasteSpecial(-4122)
$ws.Range("E235").Copy()
$ws.Range("E246").PasteSpecial(-4122)

# Row 247
$ws.Range("A247").Value = 245
$ws.Range("B247").Value = 7951775
$ws.Range("C247").Value = "Romania Liga I"
$ws.Range("D247").Value = "Romania Liga I"
$ws.Range("E247").Value = 45368.55208333334
$ws.Range("F247").Value = "AFC Hermannstadt"
$ws.Range("G247").Value = "Dinamo Bucharest"
$ws.Range("K247").Value = 1.909
$ws.Range("L247").Value = 3.1
$ws.Range("M247").Value = 4
$ws.Range("N247").Value = 1.95
$ws.Range("O247").Value = 3
$ws.Range("P247").Value = 4
$ws.Range("Q247").Value = -0.5
$ws.Range("R247").Value = 2.025
$ws.Range("S247").Value = 1.825
$ws.Range("T247").Value = 1.75
$ws.Range("U247").Value = 1.8
$ws.Range("V247").Value = 2.05
$ws.Range("W247").Value = 0
$ws.Range("X247").Value = 0
$ws.Range("Y247").Value = 0
$ws.Range("Z247").Value = 0
$ws.Range("AA247").Value = 0
$ws.Range("A235").Copy()
$ws.Range("A247").PasteSpecial(-4122)
$ws.Range("E235").Copy()
$ws.Range("E247").PasteSpecial(-4122)

# Row 248
$ws.Range("A248").Value = 246
$ws.Range("B248").Value = 7953049
$ws.Range("C248").Value = "Romania Liga I"
$ws.Range("D248").Value = "Romania Liga I"
$ws.Range("E248").Value = 45368.66666666666
$ws.Range("F248").Value = "FCSB"
$ws.Range("G248").Value = "ACS Sepsi"
$ws.Range("K248").Value = 1.65
$ws.Range("L248").Value = 3.5
$ws.Range("M248").Value = 5
$ws.Range("N248").Value = 1.75
$ws.Range("O248").Value = 3.4
$ws.Range("P248").Value = 4.5
$ws.Range("Q248").Value = -0.75
$ws.Range("R248").Value = 1.975
$ws.Range("S248").Value = 1.875
$ws.Range("T248").Value = 2.5
$ws.Range("U248").Value = 2
$ws.Range("V248").Value = 1.85
$ws.Range("W248").Value = 0
$ws.Range("X248").Value = 0
$ws.Range("Y248").Value = 0
$ws.Range("Z248").Value = 0
$ws.Range("AA248").Value = 0
$ws.Range("A235").Copy()
$ws.Range("A248").PasteSpecial(-4122)
$ws.Range("E235").Copy()
$ws.Range("E248").PasteSpecial(-4122)

# Row 249
$ws.Range("A249").Value = 247
$ws.Range("B249").Value = 7951776
$ws.Range("C249").Value = "Romania Liga I"
$ws.Range("D249").Value = "Romania Liga I"
$ws.Range("E249").Value = 45369.52083333334
$ws.Range("F249").Value = "ACS UTA Batrana Doamna"
$ws.Range("G249").Value = "FC Voluntari"
$ws.Range("K249").Value = 1.909
$ws.Range("L249").Value = 3.1
$ws.Range("M249").Value = 4
$ws.Range("N249").Value = 1.85
$ws.Range("O249").Value = 3.1
$ws.Range("P249").Value = 4.2
$ws.Range("Q249").Value = -0.5
$ws.Range("R249").Value = 1.925
$ws.Range("S249").Value = 1.925
$ws.Range("T249").Value = 2.25
$ws.Range("U249").Value = 1.95
$ws.Range("V249").Value = 1.9
$ws.Range("W249").Value = 0
$ws.Range("X249").Value = 0
$ws.Range("Y249").Value = 0
$ws.Range("Z249").Value = 0
$ws.Range("AA249").Value = 0
$ws.Range("A235").Copy()
$ws.Range("A249").PasteSpecial(-4122)
$ws.Range("E235").Copy()
$ws.Range("E249").PasteSpecial(-4122)

# Row 250
$ws.Range("A250").Value = 248
$ws.Range("B250").Value = 7951777
$ws.Range("C250").Value = "Romania Liga I"
$ws.Range("D250").Value = "Romania Liga I"
$ws.Range("E250").Value = 45369.64583333334
$ws.Range("F250").Value = "Petrolul Ploiesti"
$ws.Range("G250").Value = "FC U Craiova 1948"
$ws.Range("K250").Value = 2.25
$ws.Range("L250").Value = 3.2
$ws.Range("M250").Value = 3
$ws.Range("N250").Value = 2.15
$ws.Range("O250").Value = 3.2
$ws.Range("P250").Value = 3.2
$ws.Range("Q250").Value = -0.25
$ws.Range("R250").Value = 1.9
$ws.Range("S250").Value = 1.95
$ws.Range("T250").Value = 2.25
$ws.Range("U250").Value = 1.95
$ws.Range("V250").Value = 1.9
$ws.Range("W250").Value = 0
$ws.Range("X250").Value = 0
$ws.Range("Y250").Value = 0
$ws.Range("Z250").Value = 0
$ws.Range("AA250").Value = 0
$ws.Range("A235").Copy()
$ws.Range("A250").PasteSpecial(-4122)
$ws.Range("E235").Copy()
$ws.Range("E250").PasteSpecial(-4122)
